$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in B2, B3, B4
$ws.Range("B2").Value = 54910
$ws.Range("B3").Value = 498
$ws.Range("B4").Value = 4

# Delete row 5 entirely (shifts cells up, removes the "quartz monzonite" row)
$ws.Range("A5:B5").Delete()
